$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")
$ws.Range("E3").Value = """Global"", ""Aciclovir"""
$ws.Range("E4").Select()
